$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.541.92'
$ws.Range("E2").Value = '  +2.82%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.126.52'
$ws.Range("E3").Value = '  +1.54%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.010'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '346.98'
$ws.Range("E5").Value = '  +0.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.007'
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5255'
$ws.Range("E7").Value = '  +1.51%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4484'
$ws.Range("E8").Value = '  +0.52%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.14'
$ws.Range("E9").Value = '  +4.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09415'
$ws.Range("E10").Value = '  -0.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.184'
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.30'
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.691'
$ws.Range("E13").Value = '  +7.04%  '
$ws.Range("E14").Value = '  +3.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.116.96'
$ws.Range("E15").Value = '  +1.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '102.74'
$ws.Range("E16").Value = '  +3.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001173'
$ws.Range("E17").Value = '  +0.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.010'
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '21.58'
$ws.Range("E19").Value = '  +3.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06731'
$ws.Range("E20").Value = '  +0.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.342'
$ws.Range("E21").Value = '  +2.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.007'
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.547.16'
$ws.Range("E23").Value = '  +2.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.78'
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.343'
$ws.Range("E25").Value = '  +1.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.363.66'
$ws.Range("E26").Value = '  +0.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.26'
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.561'
$ws.Range("E28").Value = '  +1.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '163.39'
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '135.10'
$ws.Range("E30").Value = '  +1.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.164'
$ws.Range("E31").Value = '  +0.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.780'
$ws.Range("E32").Value = '  +9.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.1063'
$ws.Range("E33").Value = '  +0.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.946'
$ws.Range("E34").Value = '  +12.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.312'
$ws.Range("E35").Value = '  +1.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.972'
$ws.Range("E36").Value = '  +0.74%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.66'
$ws.Range("E37").Value = '  +5.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02657'
$ws.Range("E38").Value = '  +3.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06873'
$ws.Range("E39").Value = '  +2.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.7130'
$ws.Range("E40").Value = '  +3.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.72'
$ws.Range("E41").Value = '  +2.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2261'
$ws.Range("E42").Value = '  -1.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.329'
$ws.Range("E43").Value = '  +3.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6948'
$ws.Range("E44").Value = '  +4.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.77'
$ws.Range("E45").Value = '  +4.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.409'
$ws.Range("E46").Value = '  +4.75%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.007'
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.301'
$ws.Range("E48").Value = '  +11.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.653'
$ws.Range("E49").Value = '  +0.69%  '
$ws.Range("E50").Value = '  +1.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.234'
$ws.Range("E51").Value = '  +1.26%  '
